$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 521.63635
$ws.Range("I4").Value = 398.77777
$ws.Range("K4").Value = 398.77777
$ws.Range("M4").Value = -284.77777

$ws.Range("H28").Value = 462.5
$ws.Range("I28").Value = 356.2
$ws.Range("K28").Value = 356.2
$ws.Range("M28").Value = 128.8

$ws.Range("H32").Value = 11365819
$ws.Range("I32").Value = 18183392
$ws.Range("J32").Value = 3194.6667
$ws.Range("K32").Value = 18183392
$ws.Range("L32").Value = 3194.6667
$ws.Range("M32").Value = -18183066
$ws.Range("N32").Value = -3846.6667

$ws.Range("H70").Value = 1194.8572
$ws.Range("I70").Value = 600
$ws.Range("J70").Value = 1294
$ws.Range("K70").Value = 1800
$ws.Range("L70").Value = 3882
$ws.Range("M70").Value = -1530
$ws.Range("N70").Value = -4422

$ws.Range("H73").Value = 1194.8572
$ws.Range("I73").Value = 600
$ws.Range("J73").Value = 1294
$ws.Range("K73").Value = 1800
$ws.Range("L73").Value = 3882
$ws.Range("M73").Value = -864
$ws.Range("N73").Value = -5754

$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()

$ws.Range("H137").Value = 1833.174
$ws.Range("I137").Value = 1651.8823
$ws.Range("K137").Value = 4955.6469
$ws.Range("M137").Value = -2405.6469

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1883684.6
$ws.Range("I32").Value = 863854.9399999999
$ws.Range("K32").Value = 863854.9399999999
$ws.Range("M32").Value = -863567.9399999999

$ws.Range("H61").Value = 2677.0908
$ws.Range("I61").Value = 1741.5
$ws.Range("K61").Value = 1741.5
$ws.Range("M61").Value = -1529.5

$ws.Range("H133").Value = 107992
$ws.Range("J133").Value = 107992
$ws.Range("L133").Value = 107992
$ws.Range("N133").Value = -113052

$ws.Range("H136").Value = 2677.0908
$ws.Range("I136").Value = 1741.5
$ws.Range("K136").Value = 5224.5
$ws.Range("M136").Value = -2674.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H50").Value = 78000
$ws.Range("J50").Value = 78000
$ws.Range("L50").Value = 78000
$ws.Range("N50").Value = -79148

$ws.Range("H94").Value = 71428880
$ws.Range("I94").Value = 86956710
$ws.Range("J94").Value = 874.8
$ws.Range("K94").Value = 86956710
$ws.Range("L94").Value = 874.8
$ws.Range("M94").Value = -86956259
$ws.Range("N94").Value = -1776.8

$ws.Range("H134").Value = 2152.96
$ws.Range("I134").Value = 1329.2307
$ws.Range("K134").Value = 3987.6921
$ws.Range("M134").Value = -1452.6921

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2843744.8
$ws.Range("I31").Value = 2064.1482
$ws.Range("J31").Value = 7357002
$ws.Range("K31").Value = 2064.1482
$ws.Range("L31").Value = 7357002
$ws.Range("M31").Value = -1769.1482
$ws.Range("N31").Value = -7357592

$ws.Range("H34").Value = 2843744.8
$ws.Range("I34").Value = 2064.1482
$ws.Range("J34").Value = 7357002
$ws.Range("K34").Value = 2064.1482
$ws.Range("L34").Value = 7357002
$ws.Range("M34").Value = -1862.1482
$ws.Range("N34").Value = -7357406

$ws.Range("H58").Value = 2132.2222
$ws.Range("I58").Value = 1569.5714
$ws.Range("K58").Value = 1569.5714
$ws.Range("M58").Value = -1366.5714

$ws.Range("H69").Value = 37499.75
$ws.Range("I69").Value = 34999
$ws.Range("J69").Value = 38333.332
$ws.Range("K69").Value = 34999
$ws.Range("L69").Value = 38333.332
$ws.Range("M69").Value = -34250
$ws.Range("N69").Value = -39831.332

$ws.Range("H72").Value = 37499.75
$ws.Range("I72").Value = 34999
$ws.Range("J72").Value = 38333.332
$ws.Range("K72").Value = 104997
$ws.Range("L72").Value = 114999.996
$ws.Range("M72").Value = -101253
$ws.Range("N72").Value = -122487.996

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H107").Value = 7144114
$ws.Range("J107").Value = 1928.6666
$ws.Range("L107").Value = 1928.6666
$ws.Range("N107").Value = -5768.6666

$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

$ws.Range("H119").Value = 100000
$ws.Range("J119").Value = 100000
$ws.Range("L119").Value = 100000
$ws.Range("N119").Value = -109676

$ws.Range("H123").Value = 57582.5
$ws.Range("J123").Value = 57582.5
$ws.Range("L123").Value = 57582.5
$ws.Range("N123").Value = -67382.5

$ws.Range("H132").Value = 3664.1353
$ws.Range("I132").Value = 3064.5862
$ws.Range("K132").Value = 9193.758600000001
$ws.Range("M132").Value = -6663.758600000001

$ws.Range("H134").Value = 3754.6562
$ws.Range("I134").Value = 4195.9585
$ws.Range("K134").Value = 12587.8755
$ws.Range("M134").Value = -10052.8755

$ws.Range("H136").Value = 2132.2222
$ws.Range("I136").Value = 1569.5714
$ws.Range("K136").Value = 4708.7142
$ws.Range("M136").Value = -2158.7142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 8770.323
$ws.Range("I4").Value = 6363.526
$ws.Range("K4").Value = 19090.578
$ws.Range("M4").Value = -18978.578

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1335.6342
$ws.Range("I102").Value = 1015.7059
$ws.Range("K102").Value = 1015.7059
$ws.Range("M102").Value = 606.2941

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 21962.1
$ws.Range("I40").Value = 23779.555
$ws.Range("K40").Value = 23779.555
$ws.Range("M40").Value = -23643.555

$ws.Range("H55").Value = 549.375
$ws.Range("I55").Value = 349.5
$ws.Range("K55").Value = 349.5
$ws.Range("M55").Value = -176.5

$ws.Range("H68").Value = 1806.2858
$ws.Range("I68").Value = 1707.3334
$ws.Range("K68").Value = 1707.3334
$ws.Range("M68").Value = -958.3334

$ws.Range("H71").Value = 1806.2858
$ws.Range("I71").Value = 1707.3334
$ws.Range("K71").Value = 8536.666999999999
$ws.Range("M71").Value = -4792.666999999999

$ws.Range("H100").Value = 2614.6155
$ws.Range("I100").Value = 2855.6667
$ws.Range("J100").Value = 2072.25
$ws.Range("K100").Value = 2855.6667
$ws.Range("L100").Value = 2072.25
$ws.Range("M100").Value = -2314.6667
$ws.Range("N100").Value = -3154.25

$ws.Range("H136").Value = 7237.5835
$ws.Range("I136").Value = 7179
$ws.Range("J136").Value = 7296.1665
$ws.Range("K136").Value = 21537
$ws.Range("L136").Value = 21888.4995
$ws.Range("M136").Value = -18987
$ws.Range("N136").Value = -26988.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3211
$ws.Range("I81").Value = 4448.5386
$ws.Range("J81").Value = 1200
$ws.Range("K81").Value = 8897.0772
$ws.Range("L81").Value = 2400
$ws.Range("M81").Value = -7836.0772
$ws.Range("N81").Value = -4522

$ws.Range("H84").Value = 3211
$ws.Range("I84").Value = 4448.5386
$ws.Range("J84").Value = 1200
$ws.Range("K84").Value = 44485.386
$ws.Range("L84").Value = 12000
$ws.Range("M84").Value = -39181.386
$ws.Range("N84").Value = -22608

$ws.Range("H132").Value = 2909.2415
$ws.Range("I132").Value = 2954.8
$ws.Range("K132").Value = 8864.400000000001
$ws.Range("M132").Value = -6334.400000000001
